$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1922.615322677455
$ws.Range("D2").Value = 2496.664479588268
